$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '68.755.69'
$ws.Range('E2').Value = '  -3.75%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '3.700.48'
$ws.Range('E3').Value = '  -4.53%  '
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '0.999'
$ws.Range('E4').Value = '  -0.04%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '595.28'
$ws.Range('E5').Value = '  -0.98%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '181.33'
$ws.Range('E6').Value = '  +4.99%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '3.699.55'
$ws.Range('E7').Value = '  -4.44%  '
$ws.Range('E8').Value = '  -5.75%  '
$ws.Range('E9').Value = '  -0.03%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.718'
$ws.Range('E10').Value = '  -4.22%  '
$ws.Range('E11').Value = '  -8.89%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '56.30'
$ws.Range('E12').Value = '  +4.46%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '0.0000295'
$ws.Range('E13').Value = '  -8.59%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '10.70'
$ws.Range('E14').Value = '  -7.40%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '4.278.43'
$ws.Range('E15').Value = '  -4.79%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '3.687.80'
$ws.Range('E16').Value = '  -4.83%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '19.45'
$ws.Range('E17').Value = '  -8.16%  '
$ws.Range('E18').Value = '  -2.00%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '12.92'
$ws.Range('E19').Value = '  -7.56%  '
$ws.Range('E20').Value = '  -7.06%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '68.563.32'
$ws.Range('E21').Value = '  -3.79%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '413.49'
$ws.Range('E22').Value = '  -6.16%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '4.61'
$ws.Range('E23').Value = '  -3.39%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '89.02'
$ws.Range('E24').Value = '  -5.89%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '3.05'
$ws.Range('E25').Value = '  -7.93%  '
$ws.Range('B26').Value = 'InternetComputer(DFINITY)'
$ws.Range('C26').Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '12.74'
$ws.Range('E26').Value = '  -8.12%  '
$ws.Range('B27').Value = 'RenderToken'
$ws.Range('C27').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '10.90'
$ws.Range('E27').Value = '  -7.13%  '
$ws.Range('B28').Value = 'Toncoin'
$ws.Range('C28').Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '3.89'
$ws.Range('E28').Value = '  -3.66%  '
$ws.Range('E29').Value = '  +1.82%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '9.66'
$ws.Range('E30').Value = '  -7.76%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '33.02'
$ws.Range('E31').Value = '  -6.30%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '7.40'
$ws.Range('E32').Value = '  -14.03%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '12.50'
$ws.Range('E33').Value = '  -7.97%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '0.120'
$ws.Range('E34').Value = '  -5.54%  '
$ws.Range('B35').Value = 'InjectiveProtocol'
$ws.Range('C35').Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '43.67'
$ws.Range('E35').Value = '  -9.00%  '
$ws.Range('B36').Value = 'OKB'
$ws.Range('C36').Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '64.93'
$ws.Range('E36').Value = '  -6.81%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '601.87'
$ws.Range('E37').Value = '  -4.87%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.0₃0885'
$ws.Range('E38').Value = '  -11.78%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.405'
$ws.Range('E39').Value = '  -7.77%  '
$ws.Range('E40').Value = '  +0.26%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '1.00'
$ws.Range('E41').Value = '  -0.15%  '
$ws.Range('E42').Value = '  -6.03%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '3.06'
$ws.Range('E43').Value = '  -7.00%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.0443'
$ws.Range('E44').Value = '  -6.38%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '2.67'
$ws.Range('E45').Value = '  -6.87%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '2.78'
$ws.Range('E46').Value = '  -12.50%  '
$ws.Range('B47').Value = 'THORChain'
$ws.Range('C47').Value = 'https://coinranking.com/coin/ybmU-kKU+thorchain-rune'
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '9.12'
$ws.Range('E47').Value = '  -10.50%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '0.136'
$ws.Range('E48').Value = '  -6.19%  '
$ws.Range('B49').Value = 'WEMIXToken'
$ws.Range('C49').Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '2.73'
$ws.Range('E49').Value = '  -6.99%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '2.741.61'
$ws.Range('E50').Value = '  -5.79%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '3.09'
$ws.Range('E51').Value = '  -4.84%  '
